$d = $word.ActiveDocument

# Locate the "RUT" Heading 2 paragraph, then remove the paragraph that
# immediately follows it (the italic "路得記" title line). Deleting the
# whole paragraph Range (start..end, which includes its trailing
# paragraph mark) merges it away cleanly without disturbing the runs of
# the preceding "RUT" paragraph.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text -eq "RUT`r") {
        $target = $d.Paragraphs.Item($i + 1)
        $r = $d.Range($target.Range.Start, $target.Range.End)
        $r.Delete()
        break
    }
}
